$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")
$ws.Rows.Item(10).Delete()
$ws.Range("A10").Select() | Out-Null
